# Journal de bord CPNVoiturage - add three new journal entries (rows 60-62)
# and update the saved selection/scroll state to reflect where the author
# was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 60
$ws.Range("C60").Value = "Corrections des requêtes pour afficher les voitures"
$ws.Range("D60").Value = 44277
$ws.Range("E60").Value = 90
$ws.Range("F60").Value = "Correction d'un problème de requête qui donnait des doublons pour les conducteurs"

# Row 61
$ws.Range("C61").Value = "Changement du menu si l'on est déjà conducteur"
$ws.Range("D61").Value = 44277
$ws.Range("E61").Value = 45
$ws.Range("F61").Value = "Affichage si l'on est conducteur ou si l'on est déjà dans une voiture"

# Row 62
$ws.Range("C62").Value = "Fonction pour rejoindre une voiture"
$ws.Range("D62").Value = 44277
$ws.Range("E62").Value = 45
# F62 stays blank

# Update the saved view state: scrolled back to column A and the active
# selection moved to C62 (last entered cell).
[void]$ws.Range("C62").Select()
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
